$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("G2").Value = [double]"20.32821866666667"
$ws.Range("H2").Value = [double]"60.984656"
$ws.Range("I2").Value = [double]"0.004181898474048532"
$ws.Range("J2").Value = [double]"0.004181898474048532"
$ws.Range("K2").Value = [double]"3"
$ws.Range("M2").Value = [double]"0.6186053333333333"
$ws.Range("N2").Value = [double]"1.855816"
$ws.Range("O2").Value = [double]"0.0556943868446899"
$ws.Range("P2").Value = [double]"0.0556943868446899"
$ws.Range("Q2").Value = [double]"12.57514448436622"
$ws.Range("R2").Value = [double]"113.176300359296"
$ws.Range("S2").Value = [double]"0.0002329082713588773"
$ws.Range("T2").Value = [double]"0.0002329082713588773"
$ws.Range("E3").Value = [double]"3"
$ws.Range("G3").Value = [double]"20.32821866666667"
$ws.Range("H3").Value = [double]"60.984656"
$ws.Range("I3").Value = [double]"0.004181898474048532"
$ws.Range("J3").Value = [double]"0.004181898474048532"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"2.338622"
$ws.Range("N3").Value = [double]"7.015866"
$ws.Range("O3").Value = [double]"0.2105512373287584"
$ws.Range("P3").Value = [double]"0.2105512373287584"
$ws.Range("Q3").Value = [double]"47.54001939467734"
$ws.Range("R3").Value = [double]"427.860174552096"
$ws.Range("S3").Value = [double]"0.0008805038980941652"
$ws.Range("T3").Value = [double]"0.0008805038980941652"
$ws.Range("E4").Value = [double]"3"
$ws.Range("G4").Value = [double]"20.32821866666667"
$ws.Range("H4").Value = [double]"60.984656"
$ws.Range("I4").Value = [double]"0.004181898474048532"
$ws.Range("J4").Value = [double]"0.004181898474048532"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"0.165314"
$ws.Range("N4").Value = [double]"0.495942"
$ws.Range("O4").Value = [double]"0.01488357983794147"
$ws.Range("P4").Value = [double]"0.01488357983794148"
$ws.Range("Q4").Value = [double]"3.360539140661333"
$ws.Range("R4").Value = [double]"30.244852265952"
$ws.Range("S4").Value = [double]"6.224161981266695E-05"
$ws.Range("T4").Value = [double]"6.224161981266696E-05"
$ws.Range("E5").Value = [double]"3"
$ws.Range("G5").Value = [double]"20.32821866666667"
$ws.Range("H5").Value = [double]"60.984656"
$ws.Range("I5").Value = [double]"0.004181898474048532"
$ws.Range("J5").Value = [double]"0.004181898474048532"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"7.984598333333333"
$ws.Range("N5").Value = [double]"23.953795"
$ws.Range("O5").Value = [double]"0.7188707959886103"
$ws.Range("P5").Value = [double]"0.7188707959886103"
$ws.Range("Q5").Value = [double]"162.3126608855022"
$ws.Range("R5").Value = [double]"1460.81394796952"
$ws.Range("S5").Value = [double]"0.003006244684782823"
$ws.Range("T5").Value = [double]"0.003006244684782823"
$ws.Range("E6").Value = [double]"3"
$ws.Range("G6").Value = [double]"4809.896321333334"
$ws.Range("H6").Value = [double]"14429.688964"
$ws.Range("I6").Value = [double]"0.9894865072215304"
$ws.Range("J6").Value = [double]"0.9894865072215304"
$ws.Range("K6").Value = [double]"3"
$ws.Range("M6").Value = [double]"0.6186053333333333"
$ws.Range("N6").Value = [double]"1.855816"
$ws.Range("O6").Value = [double]"0.0556943868446899"
$ws.Range("P6").Value = [double]"0.0556943868446899"
$ws.Range("Q6").Value = [double]"2975.427517157181"
$ws.Range("R6").Value = [double]"26778.84765441463"
$ws.Range("S6").Value = [double]"0.05510884431079696"
$ws.Range("T6").Value = [double]"0.05510884431079696"
$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"4809.896321333334"
$ws.Range("H7").Value = [double]"14429.688964"
$ws.Range("I7").Value = [double]"0.9894865072215304"
$ws.Range("J7").Value = [double]"0.9894865072215304"
$ws.Range("K7").Value = [double]"3"
$ws.Range("M7").Value = [double]"2.338622"
$ws.Range("N7").Value = [double]"7.015866"
$ws.Range("O7").Value = [double]"0.2105512373287584"
$ws.Range("P7").Value = [double]"0.2105512373287584"
$ws.Range("Q7").Value = [double]"11248.5293547892"
$ws.Range("R7").Value = [double]"101236.7641931028"
$ws.Range("S7").Value = [double]"0.2083376084156047"
$ws.Range("T7").Value = [double]"0.2083376084156047"
$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"4809.896321333334"
$ws.Range("H8").Value = [double]"14429.688964"
$ws.Range("I8").Value = [double]"0.9894865072215304"
$ws.Range("J8").Value = [double]"0.9894865072215304"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"0.165314"
$ws.Range("N8").Value = [double]"0.495942"
$ws.Range("O8").Value = [double]"0.01488357983794147"
$ws.Range("P8").Value = [double]"0.01488357983794148"
$ws.Range("Q8").Value = [double]"795.1432004648987"
$ws.Range("R8").Value = [double]"7156.288804184088"
$ws.Range("S8").Value = [double]"0.0147271014287975"
$ws.Range("T8").Value = [double]"0.0147271014287975"
$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"4809.896321333334"
$ws.Range("H9").Value = [double]"14429.688964"
$ws.Range("I9").Value = [double]"0.9894865072215304"
$ws.Range("J9").Value = [double]"0.9894865072215304"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"7.984598333333333"
$ws.Range("N9").Value = [double]"23.953795"
$ws.Range("O9").Value = [double]"0.7188707959886103"
$ws.Range("P9").Value = [double]"0.7188707959886103"
$ws.Range("Q9").Value = [double]"38405.09015082427"
$ws.Range("R9").Value = [double]"345645.8113574184"
$ws.Range("S9").Value = [double]"0.7113129530663314"
$ws.Range("T9").Value = [double]"0.7113129530663314"
$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"2.69506"
$ws.Range("H10").Value = [double]"8.085180000000001"
$ws.Range("I10").Value = [double]"0.000554424737665286"
$ws.Range("J10").Value = [double]"0.000554424737665286"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"0.6186053333333333"
$ws.Range("N10").Value = [double]"1.855816"
$ws.Range("O10").Value = [double]"0.0556943868446899"
$ws.Range("P10").Value = [double]"0.0556943868446899"
$ws.Range("Q10").Value = [double]"1.667178489653333"
$ws.Range("R10").Value = [double]"15.00460640688"
$ws.Range("S10").Value = [double]"3.087834581579616E-05"
$ws.Range("T10").Value = [double]"3.087834581579616E-05"
$ws.Range("E11").Value = [double]"3"
$ws.Range("G11").Value = [double]"2.69506"
$ws.Range("H11").Value = [double]"8.085180000000001"
$ws.Range("I11").Value = [double]"0.000554424737665286"
$ws.Range("J11").Value = [double]"0.000554424737665286"
$ws.Range("K11").Value = [double]"3"
$ws.Range("M11").Value = [double]"2.338622"
$ws.Range("N11").Value = [double]"7.015866"
$ws.Range("O11").Value = [double]"0.2105512373287584"
$ws.Range("P11").Value = [double]"0.2105512373287584"
$ws.Range("Q11").Value = [double]"6.30272660732"
$ws.Range("R11").Value = [double]"56.72453946588001"
$ws.Range("S11").Value = [double]"0.0001167348145210983"
$ws.Range("T11").Value = [double]"0.0001167348145210983"
$ws.Range("E12").Value = [double]"3"
$ws.Range("G12").Value = [double]"2.69506"
$ws.Range("H12").Value = [double]"8.085180000000001"
$ws.Range("I12").Value = [double]"0.000554424737665286"
$ws.Range("J12").Value = [double]"0.000554424737665286"
$ws.Range("K12").Value = [double]"3"
$ws.Range("M12").Value = [double]"0.165314"
$ws.Range("N12").Value = [double]"0.495942"
$ws.Range("O12").Value = [double]"0.01488357983794147"
$ws.Range("P12").Value = [double]"0.01488357983794148"
$ws.Range("Q12").Value = [double]"0.44553114884"
$ws.Range("R12").Value = [double]"4.009780339560001"
$ws.Range("S12").Value = [double]"8.251824847171043E-06"
$ws.Range("T12").Value = [double]"8.251824847171043E-06"
$ws.Range("E13").Value = [double]"3"
$ws.Range("G13").Value = [double]"2.69506"
$ws.Range("H13").Value = [double]"8.085180000000001"
$ws.Range("I13").Value = [double]"0.000554424737665286"
$ws.Range("J13").Value = [double]"0.000554424737665286"
$ws.Range("K13").Value = [double]"3"
$ws.Range("M13").Value = [double]"7.984598333333333"
$ws.Range("N13").Value = [double]"23.953795"
$ws.Range("O13").Value = [double]"0.7188707959886103"
$ws.Range("P13").Value = [double]"0.7188707959886103"
$ws.Range("Q13").Value = [double]"21.51897158423333"
$ws.Range("R13").Value = [double]"193.6707442581"
$ws.Range("S13").Value = [double]"0.0003985597524812206"
$ws.Range("T13").Value = [double]"0.0003985597524812206"
$ws.Range("E14").Value = [double]"3"
$ws.Range("G14").Value = [double]"28.08283533333333"
$ws.Range("H14").Value = [double]"84.24850599999999"
$ws.Range("I14").Value = [double]"0.005777169566755752"
$ws.Range("J14").Value = [double]"0.005777169566755752"
$ws.Range("K14").Value = [double]"3"
$ws.Range("M14").Value = [double]"0.6186053333333333"
$ws.Range("N14").Value = [double]"1.855816"
$ws.Range("O14").Value = [double]"0.0556943868446899"
$ws.Range("P14").Value = [double]"0.0556943868446899"
$ws.Range("Q14").Value = [double]"17.37219171232178"
$ws.Range("R14").Value = [double]"156.349725410896"
$ws.Range("S14").Value = [double]"0.0003217559167182644"
$ws.Range("T14").Value = [double]"0.0003217559167182644"
$ws.Range("E15").Value = [double]"3"
$ws.Range("G15").Value = [double]"28.08283533333333"
$ws.Range("H15").Value = [double]"84.24850599999999"
$ws.Range("I15").Value = [double]"0.005777169566755752"
$ws.Range("J15").Value = [double]"0.005777169566755752"
$ws.Range("K15").Value = [double]"3"
$ws.Range("M15").Value = [double]"2.338622"
$ws.Range("N15").Value = [double]"7.015866"
$ws.Range("O15").Value = [double]"0.2105512373287584"
$ws.Range("P15").Value = [double]"0.2105512373287584"
$ws.Range("Q15").Value = [double]"65.67513653291066"
$ws.Range("R15").Value = [double]"591.0762287961959"
$ws.Range("S15").Value = [double]"0.001216390200538471"
$ws.Range("T15").Value = [double]"0.001216390200538471"
$ws.Range("E16").Value = [double]"3"
$ws.Range("G16").Value = [double]"28.08283533333333"
$ws.Range("H16").Value = [double]"84.24850599999999"
$ws.Range("I16").Value = [double]"0.005777169566755752"
$ws.Range("J16").Value = [double]"0.005777169566755752"
$ws.Range("K16").Value = [double]"3"
$ws.Range("M16").Value = [double]"0.165314"
$ws.Range("N16").Value = [double]"0.495942"
$ws.Range("O16").Value = [double]"0.01488357983794147"
$ws.Range("P16").Value = [double]"0.01488357983794148"
$ws.Range("Q16").Value = [double]"4.642485840294666"
$ws.Range("R16").Value = [double]"41.78237256265199"
$ws.Range("S16").Value = [double]"8.5984964484135E-05"
$ws.Range("T16").Value = [double]"8.598496448413501E-05"
$ws.Range("E17").Value = [double]"3"
$ws.Range("G17").Value = [double]"28.08283533333333"
$ws.Range("H17").Value = [double]"84.24850599999999"
$ws.Range("I17").Value = [double]"0.005777169566755752"
$ws.Range("J17").Value = [double]"0.005777169566755752"
$ws.Range("K17").Value = [double]"3"
$ws.Range("M17").Value = [double]"7.984598333333333"
$ws.Range("N17").Value = [double]"23.953795"
$ws.Range("O17").Value = [double]"0.7188707959886103"
$ws.Range("P17").Value = [double]"0.7188707959886103"
$ws.Range("Q17").Value = [double]"224.2301601978078"
$ws.Range("R17").Value = [double]"2018.07144178027"
$ws.Range("S17").Value = [double]"0.004153038485014882"
$ws.Range("T17").Value = [double]"0.004153038485014882"
